# Update scripts with new TPM data.
# Rows 5-7 (Sending cluster = MuSCs) are removed entirely, and the
# remaining rows 2-4 get refreshed numeric values from the new TPM run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete rows (previously rows 5, 6 and 7 - "MuSCs" sender).
$ws.Rows.Item(5).Resize(3).Delete()

# Row 2: FAPs -> Avp/Avpr1a -> ECs
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3143816666666667
$ws.Range("H2").Value = 0.943145
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2972546666666667
$ws.Range("N2").Value = 0.891764
$ws.Range("O2").Value = 0.02381895177797704
$ws.Range("P2").Value = 0.02381895177797704
$ws.Range("Q2").Value = 0.09345141753111111
$ws.Range("R2").Value = 0.8410627577800001
$ws.Range("S2").Value = 0.02381895177797704
$ws.Range("T2").Value = 0.02381895177797704

# Row 3: FAPs -> Avp/Avpr1a -> FAPs
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3143816666666667
$ws.Range("H3").Value = 0.943145
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.750671000000001
$ws.Range("N3").Value = 29.252013
$ws.Range("O3").Value = 0.7813191461594744
$ws.Range("P3").Value = 0.7813191461594745
$ws.Range("Q3").Value = 3.065432200098333
$ws.Range("R3").Value = 27.588889800885
$ws.Range("S3").Value = 0.7813191461594744
$ws.Range("T3").Value = 0.7813191461594745

# Row 4: FAPs -> Avp/Avpr1a -> MuSCs
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3143816666666667
$ws.Range("H4").Value = 0.943145
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.431828666666667
$ws.Range("N4").Value = 7.295485999999999
$ws.Range("O4").Value = 0.1948619020625486
$ws.Range("P4").Value = 0.1948619020625486
$ws.Range("Q4").Value = 0.7645223492744444
$ws.Range("R4").Value = 6.88070114347
$ws.Range("S4").Value = 0.1948619020625486
$ws.Range("T4").Value = 0.1948619020625486
